$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before DL (shifts DL.. onward right by one, through MN -> MO)
$ws.Range("DL1").EntireColumn.Insert()

# Set header for the newly inserted column
$ws.Range("DL1").Value = "DemonstrationProjectIdentifier"

# Update the Id value in column A for the data rows (2-7)
$ws.Range("A2:A7").Value = "6901488a7e79911955eafebd"
